$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: GRID construction - update cost-per-unit formula and units
$ws.Range("C19").Formula = "=11.5/0.5 / 1000000"
$ws.Range("E19").Value = "kilowatt"
$ws.Range("F19").Value = "Typical power size connected to the high-voltage electricity grid: P_HV = 500 MW, mean high-voltage electricity grid transportation length: l_HV = 11.5 km (Schnidrig et al., 2023)"

# Update view state to match
$ws.Range("C19").Select()
